$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing header
# cell onto the two new header cells before filling in their text.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("G2").Value = 0.125854933266722
$ws.Range("H2").Value = 0.9890000000000001
